# ===========================================================================
# Atualizacao da tabela "dados_comparados_manaus_sp"
# A serie historica avanca um dia: cada linha passa a exibir os dados do dia
# seguinte e uma nova linha (11) e adicionada com o dia mais recente coletado.
# ===========================================================================

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# As colunas de umidade (D, E, J, K) guardam valores como "61%". Se forem
# atribuidas diretamente, o Excel converte o texto para um numero percentual.
# Forcamos o formato de texto ("@") antes de escrever para preservar o "%"
# como parte do proprio texto da celula, tal como no arquivo original.
$dataRange = $ws.Range("A2:M11")
$dataRange.NumberFormat = "@"

# Linha 2: agora representa qua. 25 (antes era ter. 24)
$ws.Range("A2").Value = "qua. 25"
$ws.Range("B2").Value = "35°"
$ws.Range("C2").Value = "26°"
$ws.Range("D2").Value = "61%"
$ws.Range("E2").Value = "83%"
$ws.Range("I2").Value = "21°"
$ws.Range("J2").Value = "39%"
$ws.Range("K2").Value = "55%"

# Linha 3: agora representa qui. 26 (antes era qua. 25)
$ws.Range("A3").Value = "qui. 26"
$ws.Range("B3").Value = "35°"
$ws.Range("D3").Value = "61%"
$ws.Range("E3").Value = "82%"
$ws.Range("H3").Value = "36°"
$ws.Range("I3").Value = "22°"
$ws.Range("J3").Value = "31%"
$ws.Range("K3").Value = "48%"

# Linha 4: agora representa sex. 27 (antes era qui. 26)
$ws.Range("A4").Value = "sex. 27"
$ws.Range("B4").Value = "34°"
$ws.Range("C4").Value = "26°"
$ws.Range("D4").Value = "64%"
$ws.Range("H4").Value = "29°"
$ws.Range("I4").Value = "16°"
$ws.Range("J4").Value = "60%"
$ws.Range("K4").Value = "84%"

# Linha 5: agora representa sáb. 28 (antes era sex. 27)
$ws.Range("A5").Value = "sáb. 28"
$ws.Range("B5").Value = "34°"
$ws.Range("D5").Value = "66%"
$ws.Range("E5").Value = "86%"
$ws.Range("H5").Value = "21°"
$ws.Range("I5").Value = "15°"
$ws.Range("J5").Value = "72%"
$ws.Range("K5").Value = "85%"
$ws.Range("L5").Value = "7 de 11"

# Linha 6: agora representa dom. 29 (antes era sáb. 28)
$ws.Range("A6").Value = "dom. 29"
$ws.Range("B6").Value = "32°"
$ws.Range("C6").Value = "25°"
$ws.Range("D6").Value = "73%"
$ws.Range("E6").Value = "87%"
$ws.Range("H6").Value = "26°"
$ws.Range("J6").Value = "67%"
$ws.Range("L6").Value = "10 de 11"

# Linha 7: agora representa seg. 30 (antes era dom. 29)
$ws.Range("A7").Value = "seg. 30"
$ws.Range("B7").Value = "32°"
$ws.Range("D7").Value = "72%"
$ws.Range("E7").Value = "87%"
$ws.Range("H7").Value = "32°"
$ws.Range("I7").Value = "19°"
$ws.Range("J7").Value = "53%"
$ws.Range("K7").Value = "67%"

# Linha 8: agora representa ter. 01 (antes era seg. 30)
$ws.Range("A8").Value = "ter. 01"
$ws.Range("B8").Value = "31°"
$ws.Range("D8").Value = "75%"
$ws.Range("E8").Value = "91%"
$ws.Range("H8").Value = "35°"
$ws.Range("I8").Value = "21°"
$ws.Range("J8").Value = "36%"
$ws.Range("K8").Value = "52%"
$ws.Range("L8").Value = "Extremo"

# Linha 9: agora representa qua. 02 (antes era ter. 01)
$ws.Range("A9").Value = "qua. 02"
$ws.Range("B9").Value = "34°"
$ws.Range("D9").Value = "68%"
$ws.Range("E9").Value = "89%"
$ws.Range("H9").Value = "36°"
$ws.Range("I9").Value = "20°"
$ws.Range("K9").Value = "61%"

# Linha 10: agora representa qui. 03 (antes era qua. 02)
$ws.Range("A10").Value = "qui. 03"
$ws.Range("B10").Value = "33°"
$ws.Range("D10").Value = "71%"
$ws.Range("E10").Value = "91%"
$ws.Range("J10").Value = "48%"
$ws.Range("K10").Value = "76%"
$ws.Range("L10").Value = "Extremo"

# Linha 11: nova linha com o dia mais recente, sex. 04
$ws.Range("A11").Value = "sex. 04"
$ws.Range("B11").Value = "33°"
$ws.Range("C11").Value = "25°"
$ws.Range("D11").Value = "69%"
$ws.Range("E11").Value = "89%"
$ws.Range("F11").Value = "Extremo"
$ws.Range("G11").Value = "0 de 11"
$ws.Range("H11").Value = "29°"
$ws.Range("I11").Value = "19°"
$ws.Range("J11").Value = "58%"
$ws.Range("K11").Value = "79%"
$ws.Range("L11").Value = "Extremo"
$ws.Range("M11").Value = "0 de 11"

# Remove a formatacao de texto temporaria aplicada acima, devolvendo as
# celulas ao estilo padrao "Normal" (igual ao restante da planilha).
$dataRange.Style = "Normal"

